$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing B..E (now C..F) keep their
# formatting/content, and the old A (segment names) shifts into B, carrying
# its former "styled" look along with it.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted index column ("segments" label lives in
# the shifted-to B1 cell, matching the other headers' look).
$ws.Range("B1").Value = "segments"
$ws.Range("B1").Style = $ws.Range("C1").Style

# Index values 0..18 in column A, rows 2-20 -- styled like the old "name"
# column used to be (bold/border/center), since that look now belongs to
# the index column instead of the segment-name column.
$names = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)
    $cellA.Style = $cellB.Style
    $cellA.Value = $i
    $cellB.Style = "Normal"
}
